$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3A")

# Existing submission's Admission No (C5) was stored as text "25";
# normalize it to a real number, matching the rest of the column.
$ws.Range("C5").Value = 25

# New submission synced: 2026-02-08 19:07:39
$ws.Range("A6").Value = "2026-02-08 19:07:39"
$ws.Range("B6").Value = "Amina Abubakar Adam"
$ws.Range("C6").Value = "'47"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").Value = 8
